$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 4: "License Information" (Heading2) -> plain paragraph
# with a bold run reading "Aquifer Open Study Notes (Book Intros)"
# ------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr/>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr></w:r>' +
  '<w:r><w:rPr><w:b/><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr><w:t>Aquifer Open Study Notes (Book Intros)</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xml4)

# ------------------------------------------------------------------
# Paragraph 5: license body paragraph - rewrite the run sequence,
# drop the two hyperlinks (rId13 / rId14), keep plain text instead.
# ------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr/>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr></w:r>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr><w:t xml:space="preserve">This work is an adaptation of </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr><w:t>Tyndale Open Study Notes</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr><w:t xml:space="preserve"> &#169; 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr><w:t>Aquifer Open Study Notes</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="id_ID" w:bidi="id_ID"/></w:rPr><w:t>, was created by Mission Mutual and is also licensed under CC BY-SA 4.0.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5.Range.InsertXML($xml5)

# ------------------------------------------------------------------
# Paragraph 6: swap the "This PDF version..." sentence for the new
# multi-language adaptation blurb (pure text substitution).
# ------------------------------------------------------------------
$d.Paragraphs(6).Range.Find.Execute(
    "This PDF version is provided under the same license.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文).",
    2)
